# Update gh-pages output data (generated at 456a3b4).
#
# Refreshes the "想去人数" (want-to-go count, column F) and a couple of
# "最低票价" (min price, column G) figures that bilibili reported since the
# last scrape, across 展览/演出/本地生活, mirrors the same numbers into the
# consolidated 全部类型 sheet, and drops the "angela LIVE 2024" event from
# 全部类型 (it fell out of the upstream feed), shifting the rows below it
# up by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1_f = @{
  2  = 545
  3  = 974
  4  = 68
  7  = 1196
  8  = 958
  11 = 1055
  12 = 4044
  13 = 578
  14 = 145
  15 = 1707
  21 = 1096
  23 = 776
  24 = 667
  28 = 54
  29 = 1035
  30 = 1169
  32 = 2468
  33 = 285
  34 = 1438
  38 = 4060
}
foreach ($row in $ws1_f.Keys) {
  $ws1.Range("F$row").Value = $ws1_f[$row]
}
$ws1.Range("G2").Value = 49

# ---------------------------------------------------------------------
# Sheet 2: 演出
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2_f = @{
  6  = 199
  12 = 399
  13 = 348
  14 = 4145
  16 = 22
  17 = 29
  20 = 47
  23 = 265
  25 = 131
  27 = 239
  35 = 3
  37 = 19
  38 = 17
}
foreach ($row in $ws2_f.Keys) {
  $ws2.Range("F$row").Value = $ws2_f[$row]
}

# ---------------------------------------------------------------------
# Sheet 3: 本地生活
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3_f = @{
  4 = 1294
  6 = 461
  7 = 1043
  8 = 102
}
foreach ($row in $ws3_f.Keys) {
  $ws3.Range("F$row").Value = $ws3_f[$row]
}

# ---------------------------------------------------------------------
# Sheet 4: 全部类型 (union of all the above, own row numbering)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4_f = @{
  2  = 1294
  4  = 461
  5  = 1043
  6  = 545
  7  = 974
  8  = 68
  9  = 1196
  10 = 958
  14 = 199
  15 = 199
  17 = 1055
  19 = 4044
  20 = 578
  21 = 145
  22 = 1707
  25 = 348
  27 = 1096
  29 = 29
  31 = 776
  32 = 667
  36 = 54
  37 = 47
  39 = 265
  40 = 1035
  41 = 1169
  43 = 2468
  44 = 239
}
foreach ($row in $ws4_f.Keys) {
  $ws4.Range("F$row").Value = $ws4_f[$row]
}
$ws4.Range("G6").Value = 49

# Row 45 ("上海·「angela LIVE 2024」in SHANGHAI") no longer appears in the
# source feed. The row index column (A) is a plain 0-based sequence
# unrelated to the event data, so it must stay untouched; only B:I of the
# remaining rows shift up by one, and the now-duplicated last row is
# dropped (which also shrinks the sheet dimension from A1:I51 to A1:I50).

$ws4.Range("C45").Value = "上海·第三届Redamancy动漫游戏嘉年华"
$ws4.Range("D45").Value = "中山北路3300号4楼L4001号 环球港上海世嘉都市乐园"
$ws4.Range("E45").Value = "2024.05.01 10:00-05.03 17:00"
$ws4.Range("F45").Value = 1438
$ws4.Range("G45").Value = 60
$ws4.Range("H45").Value = "https://show.bilibili.com/platform/detail.html?id=82017"
$ws4.Range("I45").Value = "//i1.hdslb.com/bfs/openplatform/202402/UJkFbmo91708657659067.png"
$ws4.Range("B45").NumberFormat = "@"
$ws4.Range("B45").Value = "2024-05-01"

# B46 stays "2024-05-01" (same as before), only C:I move.
$ws4.Range("C46").Value = "上海·魔都野良神only"
$ws4.Range("D46").Value = "南京东路830号 第一百货"
$ws4.Range("E46").Value = "2024.05.01 10:00-05.01 17:00"
$ws4.Range("F46").Value = 471
$ws4.Range("G46").Value = 79
$ws4.Range("H46").Value = "https://show.bilibili.com/platform/detail.html?id=80321"
$ws4.Range("I46").Value = "//i2.hdslb.com/bfs/openplatform/202401/KBlb0enU1704358750268.jpeg"

$ws4.Range("C47").Value = "上海·第五十八届燃梦星辰国潮嘉年华-随机宅舞"
$ws4.Range("D47").Value = "周家嘴路3608号 宝龙旭辉广场"
$ws4.Range("E47").Value = "2024.05.02 10:20-05.03 16:30"
$ws4.Range("F47").Value = 8
$ws4.Range("G47").Value = 58
$ws4.Range("H47").Value = "https://show.bilibili.com/platform/detail.html?id=82761"
$ws4.Range("I47").Value = "//i0.hdslb.com/bfs/openplatform/202403/azEA4EM01710236719279.jpeg"
$ws4.Range("B47").NumberFormat = "@"
$ws4.Range("B47").Value = "2024-05-02"

$ws4.Range("C48").Value = "上海·钢琴诗人Pianoboy高至豪流行钢琴音乐会"
$ws4.Range("D48").Value = "南京西路1376号上海商城4层 商城剧院"
$ws4.Range("E48").Value = "2024.05.04 19:30-05.04 21:00"
$ws4.Range("F48").Value = 2
$ws4.Range("G48").Value = 126
$ws4.Range("H48").Value = "https://show.bilibili.com/platform/detail.html?id=82673"
$ws4.Range("I48").Value = "//i2.hdslb.com/bfs/openplatform/202403/MooHY44M1710149484564.jpeg"
$ws4.Range("B48").NumberFormat = "@"
$ws4.Range("B48").Value = "2024-05-04"

$ws4.Range("C49").Value = "上海·原神×崩坏×星铁only旅行盛宴2.0"
$ws4.Range("D49").Value = "西藏南路1号 上海大世界"
$ws4.Range("E49").Value = "2024.05.18 10:00-05.19 17:00"
$ws4.Range("F49").Value = 4060
$ws4.Range("G49").Value = 65
$ws4.Range("H49").Value = "https://show.bilibili.com/platform/detail.html?id=81276"
$ws4.Range("I49").Value = "//i2.hdslb.com/bfs/openplatform/202403/FtC04QSc1709635049920.jpeg"
$ws4.Range("B49").NumberFormat = "@"
$ws4.Range("B49").Value = "2024-05-18"

$ws4.Range("C50").Value = "上海·「多厨狂喜」白金交响乐团二次元交响音乐会"
$ws4.Range("D50").Value = "丁香路425号 上海东方艺术中心"
$ws4.Range("E50").Value = "2024.06.22 19:30-06.22 21:30"
$ws4.Range("F50").Value = 17
$ws4.Range("G50").Value = 99
$ws4.Range("H50").Value = "https://show.bilibili.com/platform/detail.html?id=82731"
$ws4.Range("I50").Value = "//i0.hdslb.com/bfs/openplatform/202403/K3AlF8sr1710230449280.jpeg"
$ws4.Range("B50").NumberFormat = "@"
$ws4.Range("B50").Value = "2024-06-22"

# Drop the now-duplicated trailing row (was row 51, already folded into 50).
$ws4.Rows.Item(51).Delete()
